$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalculate everything on next open (workbook.xml calcPr fullCalcOnLoad)
$wb.ForceFullCalculation = $true

# --- Populate the generated YCbCr GLCM data (A1:D4) ---
$ws.Range("A1").Value = 0.0013833096705887576
$ws.Range("B1").Value = 0.54450968530775845
$ws.Range("C1").Value = 0.99524868920505727
$ws.Range("D1").Value = 0.99930834516470568

$ws.Range("A2").Value = 0.062483666791445992
$ws.Range("B2").Value = 0.79391345877380304
$ws.Range("C2").Value = 0.6414630783814339
$ws.Range("D2").Value = 0.96875816660427705

$ws.Range("A3").Value = 0.00029458750348589481
$ws.Range("B3").Value = 0.12237632527672122
$ws.Range("C3").Value = 0.99937033284305055
$ws.Range("D3").Value = 0.99985270624825751

$ws.Range("A4").Value = 0.0024635462668391281
$ws.Range("B4").Value = 0.76201878763669328
$ws.Range("C4").Value = 0.98594506493353318
$ws.Range("D4").Value = 0.99876822686658084

# --- Column widths: A is wider than B:D, matching the authored layout ---
$ws.Columns.Item(1).ColumnWidth = 14.8
$ws.Columns.Item(2).ColumnWidth = 11.8
$ws.Columns.Item(3).ColumnWidth = 11.8
$ws.Columns.Item(4).ColumnWidth = 11.8

# --- Register the text / date-time number formats the workbook ends up
#     carrying (numFmtId 49 "@" and numFmtId 22 "m/d/yy h:mm") on scratch
#     cells, then clear and remove them so the final grid stays unstyled
#     while the style table keeps the two extra cellXfs entries. ---
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "m/d/yy h:mm"
$ws.Range("F1:F2").ClearFormats()
$ws.Range("F1:F2").Delete()
